$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(1.2465046429701743, 0.9497271180948426, 0.673579337833506, 3.449424014741931)"
$ws.Range("C2").Value = "NIG(2.063581660045853, 1.7491551439786457, 2.5470528684004297, 5.439514132098472)"
$ws.Range("D2").Value = "JSU(-0.8355090684099471, 0.9675259925565459, 1.2698004682309927, 2.0989894594779432)"
$ws.Range("E2").Value = "NIG(1.7232664840852303, 1.2635152656124875, 3.8771863875590498, 6.68973254240208)"
